$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 161.279784
$ws.Range("H2").Value = 483.839352
$ws.Range("I2").Value = 0.3023989599621841
$ws.Range("J2").Value = 0.3023989599621841
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 18146.13555886733
$ws.Range("R2").Value = 163315.2200298059
$ws.Range("S2").Value = 0.09904512379068982
$ws.Range("T2").Value = 0.09904512379068982

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 161.279784
$ws.Range("H3").Value = 483.839352
$ws.Range("I3").Value = 0.3023989599621841
$ws.Range("J3").Value = 0.3023989599621841
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 17146.37411255534
$ws.Range("R3").Value = 154317.3670129981
$ws.Range("S3").Value = 0.09358823210761509
$ws.Range("T3").Value = 0.0935882321076151

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 161.279784
$ws.Range("H4").Value = 483.839352
$ws.Range("I4").Value = 0.3023989599621841
$ws.Range("J4").Value = 0.3023989599621841
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 20110.24323876245
$ws.Range("R4").Value = 180992.189148862
$ws.Range("S4").Value = 0.1097656040638792
$ws.Range("T4").Value = 0.1097656040638792

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 288.7700093333333
$ws.Range("H5").Value = 866.3100279999999
$ws.Range("I5").Value = 0.541442630470476
$ws.Range("J5").Value = 0.5414426304704759
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 32490.49325796499
$ws.Range("R5").Value = 292414.4393216849
$ws.Range("S5").Value = 0.1773394074080522
$ws.Range("T5").Value = 0.1773394074080522

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 288.7700093333333
$ws.Range("H6").Value = 866.3100279999999
$ws.Range("I6").Value = 0.541442630470476
$ws.Range("J6").Value = 0.5414426304704759
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 30700.42933908834
$ws.Range("R6").Value = 276303.8640517951
$ws.Range("S6").Value = 0.1675688917044071
$ws.Range("T6").Value = 0.1675688917044071

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 288.7700093333333
$ws.Range("H7").Value = 866.3100279999999
$ws.Range("I7").Value = 0.541442630470476
$ws.Range("J7").Value = 0.5414426304704759
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 36007.21047439545
$ws.Range("R7").Value = 324064.894269559
$ws.Range("S7").Value = 0.1965343313580167
$ws.Range("T7").Value = 0.1965343313580167

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 83.28466000000002
$ws.Range("H8").Value = 249.85398
$ws.Range("I8").Value = 0.15615840956734
$ws.Range("J8").Value = 0.15615840956734
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 9370.639598166721
$ws.Range("R8").Value = 84335.75638350048
$ws.Range("S8").Value = 0.05114676653811438
$ws.Range("T8").Value = 0.05114676653811438

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 83.28466000000002
$ws.Range("H9").Value = 249.85398
$ws.Range("I9").Value = 0.15615840956734
$ws.Range("J9").Value = 0.15615840956734
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 8854.364153891562
$ws.Range("R9").Value = 79689.27738502405
$ws.Range("S9").Value = 0.04832883513214408
$ws.Range("T9").Value = 0.04832883513214409

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 83.28466000000002
$ws.Range("H10").Value = 249.85398
$ws.Range("I10").Value = 0.15615840956734
$ws.Range("J10").Value = 0.15615840956734
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 10384.90211927386
$ws.Range("R10").Value = 93464.1190734647
$ws.Range("S10").Value = 0.05668280789708149
$ws.Range("T10").Value = 0.05668280789708149
